# Aggiunta di due nuove registrazioni di attività (righe 6 e 7) al foglio.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Riga 6: Giovanni - Interno - 135 min - 29/01/2019
$ws.Range("A6").Value = "Giovanni"
$ws.Range("B6").Value = "Interno"
$ws.Range("C6").Value = 135
$ws.Range("D6").Value = 43494

# Riga 7: Giovanni - GDPR - 53 min - 29/01/2019
$ws.Range("A7").Value = "Giovanni"
$ws.Range("B7").Value = "GDPR"
$ws.Range("C7").Value = 53
$ws.Range("D7").Value = 43494

# Applica alle nuove celle data lo stesso formato delle celle data esistenti,
# riutilizzando lo stile già presente invece di crearne uno nuovo.
$ws.Range("D5").Copy()
$ws.Range("D6:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Ripristina la selezione dell'utente al termine della modifica.
$ws.Range("F11").Select()
